$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.01
$ws.Range("G2").Value = 1.01
$ws.Range("H2").Value = 1000
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 150
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 0
$ws.Range("V2").Value = 1.01
$ws.Range("W2").Value = 230
$ws.Range("X2").Value = 1000
$ws.Range("Y2").Value = 1000
$ws.Range("Z2").Value = 1000
$ws.Range("AA2").Value = 1000
$ws.Range("AB2").Value = 1000
$ws.Range("AC2").Value = 1000
$ws.Range("AD2").Value = 1000
$ws.Range("AE2").Value = 1000
$ws.Range("AF2").Value = 1000
$ws.Range("AG2").Value = 1000
$ws.Range("AH2").Value = 1000
$ws.Range("AI2").Value = 1000
$ws.Range("AJ2").Value = 1000
$ws.Range("AK2").Value = 1.13
$ws.Range("AL2").Value = 1000
$ws.Range("AM2").Value = 1000
$ws.Range("AN2").Value = 1000
$ws.Range("AO2").Value = 1000

# Row 3
$ws.Range("F3").Value = 1000
$ws.Range("G3").Value = 1000
$ws.Range("H3").Value = 1.01
$ws.Range("I3").Value = 1.01
$ws.Range("J3").Value = 1000
$ws.Range("K3").Value = 1000
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 0
$ws.Range("V3").Value = 500
$ws.Range("W3").Value = 1.01
$ws.Range("X3").Value = 1000
$ws.Range("Y3").Value = 1000
$ws.Range("Z3").Value = 1000
$ws.Range("AA3").Value = 1000
$ws.Range("AB3").Value = 1000
$ws.Range("AC3").Value = 1000
$ws.Range("AD3").Value = 1000
$ws.Range("AE3").Value = 1.21
$ws.Range("AF3").Value = 1000
$ws.Range("AG3").Value = 1000
$ws.Range("AH3").Value = 1000
$ws.Range("AI3").Value = 1000
$ws.Range("AJ3").Value = 1000
$ws.Range("AK3").Value = 1000
$ws.Range("AL3").Value = 1000
$ws.Range("AM3").Value = 1000
$ws.Range("AN3").Value = 1000
$ws.Range("AO3").Value = 1000

# Row 4
$ws.Range("F4").Value = 1.03
$ws.Range("G4").Value = 1.04
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 590
$ws.Range("J4").Value = 34
$ws.Range("K4").Value = 40
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 1.81
$ws.Range("Q4").Value = 2.14
$ws.Range("R4").Value = 1.11
$ws.Range("S4").Value = 8.8
$ws.Range("T4").Value = 3.6
$ws.Range("U4").Value = 1.32
$ws.Range("V4").Value = 1.01
$ws.Range("W4").Value = 34
$ws.Range("X4").Value = 1000
$ws.Range("Y4").Value = 1000
$ws.Range("Z4").Value = 1000
$ws.Range("AA4").Value = 1000
$ws.Range("AB4").Value = 1000
$ws.Range("AC4").Value = 1000
$ws.Range("AD4").Value = 1000
$ws.Range("AE4").Value = 1000
$ws.Range("AF4").Value = 1.88
$ws.Range("AG4").Value = 6.6
$ws.Range("AH4").Value = 55
$ws.Range("AI4").Value = 1000
$ws.Range("AJ4").Value = 5.8
$ws.Range("AK4").Value = 23
$ws.Range("AL4").Value = 220
$ws.Range("AM4").Value = 1000
$ws.Range("AN4").Value = 32
$ws.Range("AO4").Value = 1000

# Row 5
$ws.Range("F5").Value = 3.7
$ws.Range("G5").Value = 4
$ws.Range("H5").Value = 6.2
$ws.Range("I5").Value = 6.8
$ws.Range("J5").Value = 1.7
$ws.Range("K5").Value = 1.75
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 0
$ws.Range("T5").Value = 0
$ws.Range("U5").Value = 0
$ws.Range("V5").Value = 1.17
$ws.Range("W5").Value = 1.3
$ws.Range("X5").Value = 1000
$ws.Range("Y5").Value = 1000
$ws.Range("Z5").Value = 1000
$ws.Range("AA5").Value = 1000
$ws.Range("AB5").Value = 1000
$ws.Range("AC5").Value = 1000
$ws.Range("AD5").Value = 1000
$ws.Range("AE5").Value = 1000
$ws.Range("AF5").Value = 1000
$ws.Range("AG5").Value = 1000
$ws.Range("AH5").Value = 1.84
$ws.Range("AI5").Value = 7.8
$ws.Range("AJ5").Value = 1000
$ws.Range("AK5").Value = 1000
$ws.Range("AL5").Value = 6
$ws.Range("AM5").Value = 28
$ws.Range("AN5").Value = 36
$ws.Range("AO5").Value = 65

# Row 6
$ws.Range("F6").Value = 1.75
$ws.Range("G6").Value = 1.81
$ws.Range("J6").Value = 3.85
$ws.Range("K6").Value = 4
$ws.Range("M6").Value = 1.09
$ws.Range("N6").Value = 3.2
$ws.Range("P6").Value = 1.75
$ws.Range("Q6").Value = 2.24
$ws.Range("S6").Value = 4.1
$ws.Range("T6").Value = 2.02
$ws.Range("U6").Value = 1.84
$ws.Range("V6").Value = 1.21
$ws.Range("W6").Value = 2.22
$ws.Range("X6").Value = 12.5
$ws.Range("Y6").Value = 16
$ws.Range("AE6").Value = 130
$ws.Range("AG6").Value = 10.5
$ws.Range("AI6").Value = 110
$ws.Range("AK6").Value = 22
$ws.Range("AL6").Value = 46
$ws.Range("AM6").Value = 1000
$ws.Range("AN6").Value = 16
